$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.82
$ws.Range("G2").Value = 2.88
$ws.Range("H2").Value = 2.7
$ws.Range("I2").Value = 2.82
$ws.Range("J2").Value = 3.35
$ws.Range("K2").Value = 3.55
$ws.Range("L2").Value = 1.32
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 5.1
$ws.Range("O2").Value = 1.24
$ws.Range("P2").Value = 2.28
$ws.Range("Q2").Value = 1.74
$ws.Range("R2").Value = 1.53
$ws.Range("S2").Value = 2.76
$ws.Range("T2").Value = 1.61
$ws.Range("U2").Value = 2.6
$ws.Range("V2").Value = 1.55
$ws.Range("X2").Value = 19.5
$ws.Range("Y2").Value = 24
$ws.Range("Z2").Value = 20
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 24
$ws.Range("AC2").Value = 8.4
$ws.Range("AD2").Value = 980
$ws.Range("AE2").Value = 30
$ws.Range("AF2").Value = 20
$ws.Range("AG2").Value = 12.5
$ws.Range("AH2").Value = 16
$ws.Range("AI2").Value = 40
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 980
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 960
$ws.Range("AO2").Value = 1000
$ws.Range("F3").Value = 3.85
$ws.Range("G3").Value = 4.4
$ws.Range("H3").Value = 1.79
$ws.Range("I3").Value = 1.9
$ws.Range("J3").Value = 4.4
$ws.Range("L3").Value = 1.29
$ws.Range("N3").Value = 5.7
$ws.Range("O3").Value = 1.18
$ws.Range("P3").Value = 2.64
$ws.Range("Q3").Value = 1.49
$ws.Range("R3").Value = 1.64
$ws.Range("S3").Value = 2.26
$ws.Range("T3").Value = 1.55
$ws.Range("U3").Value = 2.48
$ws.Range("V3").Value = 2.12
$ws.Range("W3").Value = 1.3
$ws.Range("X3").Value = 32
$ws.Range("Y3").Value = 14.5
$ws.Range("Z3").Value = 15
$ws.Range("AA3").Value = 21
$ws.Range("AC3").Value = 12.5
$ws.Range("AD3").Value = 11.5
$ws.Range("AH3").Value = 20
$ws.Range("AI3").Value = 42
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 44
$ws.Range("AL3").Value = 65
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("F4").Value = 11
$ws.Range("I4").Value = 1.44
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 1.36
$ws.Range("N4").Value = 3.85
$ws.Range("P4").Value = 1.95
$ws.Range("Q4").Value = 1.87
$ws.Range("U4").Value = 1.73
$ws.Range("V4").Value = 3.25
$ws.Range("W4").Value = 1.09
$ws.Range("Y4").Value = 32
$ws.Range("AE4").Value = 36
$ws.Range("AO4").Value = 15
$ws.Range("F5").Value = 1.71
$ws.Range("G5").Value = 1.8
$ws.Range("H5").Value = 5.9
$ws.Range("J5").Value = 3.5
$ws.Range("L5").Value = 1.52
$ws.Range("N5").Value = 2.7
$ws.Range("P5").Value = 1.62
$ws.Range("Q5").Value = 2.34
$ws.Range("U5").Value = 1.69
$ws.Range("V5").Value = 1.16
$ws.Range("W5").Value = 2.24
$ws.Range("AB5").Value = 11.5
$ws.Range("AF5").Value = 21
$ws.Range("J6").Value = 4.4
$ws.Range("L6").Value = 1.41
$ws.Range("U6").Value = 1.71
$ws.Range("AN6").Value = 390
$ws.Range("F7").Value = 2.88
$ws.Range("K7").Value = 3.35
$ws.Range("M7").Value = 1.09
$ws.Range("N7").Value = 3.1
$ws.Range("S7").Value = 4.2
$ws.Range("U7").Value = 2
$ws.Range("AH7").Value = 38
$ws.Range("AJ7").Value = 55
$ws.Range("U8").Value = 2.02
$ws.Range("H9").Value = 13.5
$ws.Range("J9").Value = 7.4
$ws.Range("N9").Value = 8
$ws.Range("Q9").Value = 1.37
$ws.Range("R9").Value = 1.96
$ws.Range("T9").Value = 1.83
$ws.Range("U9").Value = 1.93
$ws.Range("Y9").Value = 470
$ws.Range("AB9").Value = 16
$ws.Range("AD9").Value = 60
$ws.Range("AF9").Value = 11
$ws.Range("AG9").Value = 13.5
$ws.Range("AJ9").Value = 11
$ws.Range("AN9").Value = 3.15
$ws.Range("K10").Value = 3.35
$ws.Range("N10").Value = 2.82
$ws.Range("P10").Value = 1.62
$ws.Range("U10").Value = 1.76
$ws.Range("V10").Value = 1.24
$ws.Range("Z10").Value = 38
$ws.Range("AB10").Value = 6.8
$ws.Range("AH10").Value = 26
$ws.Range("AO10").Value = 150
$ws.Range("G11").Value = 1.23
$ws.Range("K11").Value = 9.199999999999999
$ws.Range("L11").Value = 1.16
$ws.Range("U11").Value = 2.2
$ws.Range("V11").Value = 1.05
$ws.Range("W11").Value = 5.3
$ws.Range("X11").Value = 380
$ws.Range("AA11").Value = 620
$ws.Range("AE11").Value = 220
$ws.Range("AK11").Value = 17
$ws.Range("F12").Value = 2.78
$ws.Range("P12").Value = 2.92
$ws.Range("R12").Value = 1.8
$ws.Range("T12").Value = 1.51
$ws.Range("AC12").Value = 11.5
$ws.Range("AF12").Value = 75
$ws.Range("H13").Value = 7
$ws.Range("I13").Value = 7.2
$ws.Range("N13").Value = 6
$ws.Range("O13").Value = 1.18
$ws.Range("R13").Value = 1.67
$ws.Range("S13").Value = 2.46
$ws.Range("U13").Value = 2.28
$ws.Range("Y13").Value = 32
$ws.Range("AF13").Value = 10
$ws.Range("J14").Value = 3.8
$ws.Range("T14").Value = 1.54
$ws.Range("F15").Value = 2.14
$ws.Range("G15").Value = 2.18
$ws.Range("H15").Value = 3.95
$ws.Range("I15").Value = 4.2
$ws.Range("K15").Value = 3.5
$ws.Range("R15").Value = 1.3
$ws.Range("S15").Value = 4.1
$ws.Range("T15").Value = 1.94
$ws.Range("V15").Value = 1.31
$ws.Range("W15").Value = 1.84
$ws.Range("Z15").Value = 27
$ws.Range("AA15").Value = 85
$ws.Range("AH15").Value = 21
$ws.Range("AI15").Value = 70
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = 4.3
$ws.Range("I16").Value = 4.5
$ws.Range("J16").Value = 3.55
$ws.Range("K16").Value = 3.75
$ws.Range("V16").Value = 1.28
$ws.Range("K17").Value = 16
$ws.Range("L17").Value = 1.15
$ws.Range("S17").Value = 1.58
$ws.Range("T17").Value = 1.74
$ws.Range("F18").Value = 1.85
$ws.Range("I18").Value = 5.2
$ws.Range("J18").Value = 3.4
$ws.Range("O18").Value = 1.38
$ws.Range("P18").Value = 1.74
$ws.Range("Q18").Value = 1.98
$ws.Range("R18").Value = 1.27
$ws.Range("T18").Value = 1.92
$ws.Range("U18").Value = 1.89
$ws.Range("X18").Value = 15
$ws.Range("AF18").Value = 13.5
$ws.Range("AE3").Value = 17.5
$ws.Range("AF3").Value = 65
$ws.Range("AG3").Value = 18.5
